$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.289.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = "'2.176.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'251.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.72%  '
$ws.Range("D6").Value = "'0.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").Value = "'72.60"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.77%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -3.26%  '
$ws.Range("D10").Value = "'39.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.70%  '
$ws.Range("D11").Value = "'0.0910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.59%  '
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("D13").Value = "'6.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.85%  '
$ws.Range("D14").Value = "'2.502.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'14.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.15%  '
$ws.Range("D16").Value = "'2.169.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.01%  '
$ws.Range("D17").Value = "'0.762"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.35%  '
$ws.Range("D18").Value = "'42.226.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("E19").Value = '  -3.56%  '
$ws.Range("D20").Value = "'70.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("D21").Value = "'5.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").Value = "'225.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'9.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.60%  '
$ws.Range("E24").Value = '  -1.17%  '
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("D26").Value = "'10.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("E28").Value = '  +5.96%  '
$ws.Range("D29").Value = "'2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.49%  '
$ws.Range("D30").Value = "'170.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.13%  '
$ws.Range("D31").Value = "'36.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.05%  '
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = "'0.0806"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.34%  '
$ws.Range("E34").Value = '  -5.75%  '
$ws.Range("E35").Value = '  -1.53%  '
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("E37").Value = '  -4.89%  '
$ws.Range("D38").Value = "'0.0333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.13%  '
$ws.Range("E39").Value = '  -4.45%  '
$ws.Range("D40").Value = "'11.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.31%  '
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("D42").Value = "'58.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.16%  '
$ws.Range("E43").Value = '  -7.01%  '
$ws.Range("D44").Value = "'100.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("E45").Value = '  +8.45%  '
$ws.Range("D46").Value = "'0.0970"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'8.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.58%  '
$ws.Range("B48").Value = 'WOONetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D48").Value = "'0.458"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.67%  '
$ws.Range("E49").Value = '  -2.21%  '
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("E51").Value = '  +0.09%  '
